$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the D and E columns for affected rows so that numeric-looking
# strings like "1.005" or "26.856.41" are stored as text, matching the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.856.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.825.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.71"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4569"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3678"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07153"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8726"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07764"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.55"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.823.30"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.310"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.374"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.84"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008700"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.892.70"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.45"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.995"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.056.23"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.999"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.36"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.12"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.950"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.57"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.900"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08786"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.036"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7485"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.475"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.131"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.536"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.60%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01935"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.913"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05119"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.914"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4961"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1592"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.288"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4681"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.005"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.06"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.41"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.608"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06095"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.39"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.38%  "
